# Add a "Generation Charge" column to the "Historical GC" sheet, between
# the existing "Average Generation Cost" (C) and "Date" (D) columns.
# This shifts the old D (Date) -> E and old E (Power Supplier ID) -> F.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Historical GC")

# Insert a new blank column at D; existing D/E columns (and their
# formatting) shift right to E/F.
$ws.Columns.Item(4).Insert()

# Header for the newly inserted column.
$ws.Cells.Item(1, 4).Value = "Generation Charge"

# Generation Charge values, grouped by the Date column (rows share the
# same value within each month block).
$ws.Cells.Item(2, 4).Value = "2 .7423"
$ws.Cells.Item(3, 4).Value = "2 .7423"
$ws.Cells.Item(4, 4).Value = "2 .7423"
$ws.Cells.Item(5, 4).Value = "2 .7423"

$ws.Cells.Item(6, 4).Value = "6 .7472"
$ws.Cells.Item(7, 4).Value = "6 .7472"
$ws.Cells.Item(8, 4).Value = "6 .7472"
$ws.Cells.Item(9, 4).Value = "6 .7472"
$ws.Cells.Item(10, 4).Value = "6 .7472"

$ws.Cells.Item(11, 4).Value = "4 .9881"
$ws.Cells.Item(12, 4).Value = "4 .9881"
$ws.Cells.Item(13, 4).Value = "4 .9881"
$ws.Cells.Item(14, 4).Value = "4 .9881"
$ws.Cells.Item(15, 4).Value = "4 .9881"

$ws.Cells.Item(16, 4).Value = "7 .9383"
$ws.Cells.Item(17, 4).Value = "7 .9383"
$ws.Cells.Item(18, 4).Value = "7 .9383"
$ws.Cells.Item(19, 4).Value = "7 .9383"
